$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44208
$ws.Range("J2").Value = 85
$ws.Range("K2").Value = 3700
$ws.Range("L2").Value = 4000
$ws.Range("M2").Value = 3824
$ws.Range("N2").Value = '$/paquete 2 kilos'
$ws.Range("O2").Value = 'Provincia de Diguillín'
$ws.Range("P2").Value = 1912
$ws.Range("Q2").Value = 2

# Row 3
$ws.Range("D3").Value = 44225
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 3400
$ws.Range("L3").Value = 3700
$ws.Range("M3").Value = 3550
$ws.Range("N3").Value = '$/paquete 2 kilos'
$ws.Range("O3").Value = 'Provincia de Diguillín'
$ws.Range("P3").Value = 1775
$ws.Range("Q3").Value = 2

# Row 4
$ws.Range("D4").Value = 44210
$ws.Range("J4").Value = 105
$ws.Range("K4").Value = 3500
$ws.Range("L4").Value = 4000
$ws.Range("M4").Value = 3714
$ws.Range("N4").Value = '$/paquete 2 kilos'
$ws.Range("O4").Value = 'Provincia de Diguillín'
$ws.Range("P4").Value = 1857
$ws.Range("Q4").Value = 2

# Row 5
$ws.Range("D5").Value = 44215
$ws.Range("J5").Value = 140
$ws.Range("K5").Value = 3500
$ws.Range("L5").Value = 4000
$ws.Range("M5").Value = 3768
$ws.Range("N5").Value = '$/paquete 2 kilos'
$ws.Range("O5").Value = 'Provincia de Diguillín'
$ws.Range("P5").Value = 1884
$ws.Range("Q5").Value = 2

# Row 6
$ws.Range("D6").Value = 44161
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 2800
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = 2900
$ws.Range("N6").Value = '$/paquete 2 kilos'
$ws.Range("O6").Value = 'Provincia de Diguillín'
$ws.Range("P6").Value = 1450
$ws.Range("Q6").Value = 2

# Row 7
$ws.Range("D7").Value = 44762
$ws.Range("J7").Value = 60
$ws.Range("K7").Value = 8000
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = 8000
$ws.Range("N7").Value = '$/docena de atados'
$ws.Range("O7").Value = 'Región Metropolitana'
$ws.Range("P7").Value = 2667
$ws.Range("Q7").Value = 3

# Row 8
$ws.Range("D8").Value = 44166
$ws.Range("J8").Value = 70
$ws.Range("K8").Value = 3500
$ws.Range("L8").Value = 4000
$ws.Range("M8").Value = 3679
$ws.Range("N8").Value = '$/paquete 36 unidades'
$ws.Range("O8").Value = 'Región Metropolitana'
$ws.Range("P8").Value = 102
$ws.Range("Q8").Value = 36

# Row 9
$ws.Range("D9").Value = 44760
$ws.Range("J9").Value = 120
$ws.Range("K9").Value = 8000
$ws.Range("L9").Value = 8000
$ws.Range("M9").Value = 8000
$ws.Range("N9").Value = '$/docena de atados'
$ws.Range("O9").Value = 'Región Metropolitana'
$ws.Range("P9").Value = 2667
$ws.Range("Q9").Value = 3

# Row 10
$ws.Range("D10").Value = 44769
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 8000
$ws.Range("L10").Value = 8000
$ws.Range("M10").Value = 8000
$ws.Range("N10").Value = '$/docena de atados'
$ws.Range("O10").Value = 'Provincia de Diguillín'
$ws.Range("P10").Value = 2667
$ws.Range("Q10").Value = 3

# Row 11
$ws.Range("D11").Value = 44704
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 6000
$ws.Range("L11").Value = 6500
$ws.Range("M11").Value = 6250
$ws.Range("N11").Value = '$/paquete 36 unidades'
$ws.Range("O11").Value = 'Región Metropolitana'
$ws.Range("P11").Value = 174
$ws.Range("Q11").Value = 36

# Row 12
$ws.Range("D12").Value = 44662
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 8000
$ws.Range("L12").Value = 8500
$ws.Range("M12").Value = 8250
$ws.Range("N12").Value = '$/paquete 36 unidades'
$ws.Range("O12").Value = 'Región Metropolitana'
$ws.Range("P12").Value = 229
$ws.Range("Q12").Value = 36

# Row 13
$ws.Range("D13").Value = 44209
$ws.Range("J13").Value = 150
$ws.Range("K13").Value = 3500
$ws.Range("L13").Value = 4000
$ws.Range("M13").Value = 3767
$ws.Range("N13").Value = '$/paquete 2 kilos'
$ws.Range("O13").Value = 'Provincia de Diguillín'
$ws.Range("P13").Value = 1884
$ws.Range("Q13").Value = 2

# Row 14
$ws.Range("D14").Value = 44160
$ws.Range("J14").Value = 43
$ws.Range("K14").Value = 3500
$ws.Range("L14").Value = 4000
$ws.Range("M14").Value = 3709
$ws.Range("N14").Value = '$/paquete 36 unidades'
$ws.Range("O14").Value = 'Región Metropolitana'
$ws.Range("P14").Value = 103
$ws.Range("Q14").Value = 36

# Row 15
$ws.Range("D15").Value = 44784
$ws.Range("J15").Value = 160
$ws.Range("K15").Value = 8000
$ws.Range("L15").Value = 8500
$ws.Range("M15").Value = 8250
$ws.Range("N15").Value = '$/docena de atados'
$ws.Range("O15").Value = 'Provincia de Diguillín'
$ws.Range("P15").Value = 2750
$ws.Range("Q15").Value = 3

# Row 16
$ws.Range("D16").Value = 44664
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 8000
$ws.Range("L16").Value = 8500
$ws.Range("M16").Value = 8250
$ws.Range("N16").Value = '$/paquete 36 unidades'
$ws.Range("O16").Value = 'Región Metropolitana'
$ws.Range("P16").Value = 229
$ws.Range("Q16").Value = 36

# Row 17
$ws.Range("D17").Value = 44775
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 8000
$ws.Range("L17").Value = 8000
$ws.Range("M17").Value = 8000
$ws.Range("N17").Value = '$/docena de atados'
$ws.Range("O17").Value = 'Provincia de Diguillín'
$ws.Range("P17").Value = 2667
$ws.Range("Q17").Value = 3

# Row 18
$ws.Range("D18").Value = 44764
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 8000
$ws.Range("L18").Value = 9000
$ws.Range("M18").Value = 8500
$ws.Range("N18").Value = '$/docena de atados'
$ws.Range("O18").Value = 'Región Metropolitana'
$ws.Range("P18").Value = 2833
$ws.Range("Q18").Value = 3

# Row 19
$ws.Range("D19").Value = 44223
$ws.Range("J19").Value = 80
$ws.Range("K19").Value = 3500
$ws.Range("L19").Value = 3800
$ws.Range("M19").Value = 3688
$ws.Range("N19").Value = '$/paquete 2 kilos'
$ws.Range("O19").Value = 'Provincia de Diguillín'
$ws.Range("P19").Value = 1844
$ws.Range("Q19").Value = 2

# Row 20
$ws.Range("D20").Value = 44771
$ws.Range("J20").Value = 150
$ws.Range("K20").Value = 8000
$ws.Range("L20").Value = 8000
$ws.Range("M20").Value = 8000
$ws.Range("N20").Value = '$/docena de atados'
$ws.Range("O20").Value = 'Provincia de Diguillín'
$ws.Range("P20").Value = 2667
$ws.Range("Q20").Value = 3

# Row 21
$ws.Range("D21").Value = 44782
$ws.Range("J21").Value = 100
$ws.Range("K21").Value = 8000
$ws.Range("L21").Value = 8000
$ws.Range("M21").Value = 8000
$ws.Range("N21").Value = '$/docena de atados'
$ws.Range("O21").Value = 'Provincia de Diguillín'
$ws.Range("P21").Value = 2667
$ws.Range("Q21").Value = 3

# Row 22
$ws.Range("D22").Value = 44701
$ws.Range("J22").Value = 120
$ws.Range("K22").Value = 7000
$ws.Range("L22").Value = 7500
$ws.Range("M22").Value = 7250
$ws.Range("N22").Value = '$/paquete 36 unidades'
$ws.Range("O22").Value = 'Región Metropolitana'
$ws.Range("P22").Value = 2667
$ws.Range("Q22").Value = 36
